$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the sheet ("DNB Mastercard Demo" -> "Sheet")
$ws.Name = "Sheet"

# 2) Update existing rows 2-14 in place (dates, merchants, amounts) and
#    insert the new rows 15-20 with the freshly diversified merchant data.
#    Column layout: A=Dato, B=Belopet gjelder, E=Inn, F=Ut

$ws.Range("A2").Value = 45838
$ws.Range("B2").Value = "SPOTIFY AB"
$ws.Range("F2").Value = 129

$ws.Range("A3").Value = 45837
$ws.Range("B3").Value = "KIWI LØKKA"
$ws.Range("F3").Value = 423.45

$ws.Range("A4").Value = 45836
$ws.Range("B4").Value = "NETFLIX.COM"
$ws.Range("F4").Value = 179

$ws.Range("A5").Value = 45835
$ws.Range("B5").Value = "SPORT OUTLET ALNA"
$ws.Range("F5").Value = 899

$ws.Range("A6").Value = 45833
$ws.Range("B6").Value = "REMA 1000 STORO"
$ws.Range("F6").Value = 612.8

$ws.Range("A7").Value = 45832
$ws.Range("B7").Value = "PRINCESS GRUNERLØKKA"
$ws.Range("F7").Value = 567

$ws.Range("A8").Value = 45830
$ws.Range("B8").Value = "VINMONOPOLET MAJORSTUEN"
$ws.Range("F8").Value = 534

$ws.Range("A9").Value = 45828
$ws.Range("B9").Value = "STARBUCKS GRØNLAND"
$ws.Range("F9").Value = 79

$ws.Range("A10").Value = 45826
$ws.Range("B10").Value = "JERNIA STORO"
$ws.Range("F10").Value = 345

$ws.Range("A11").Value = 45825
$ws.Range("B11").Value = "BOHUS ALNA"
$ws.Range("F11").Value = 4567

$ws.Range("A12").Value = 45823
$ws.Range("B12").Value = "SKEIDAR STORO"
$ws.Range("F12").Value = 2345

# Row 13 used to be the "Innbetaling" (deposit, column E) row; it's now a
# regular merchant spend row in column F, so clear out the old E13 value.
$ws.Range("A13").Value = 45822
$ws.Range("B13").Value = "TOYS R US OSLO"
$ws.Range("E13").ClearContents()
$ws.Range("F13").Value = 567

$ws.Range("A14").Value = 45820
$ws.Range("B14").Value = "MENY CC VEST"
$ws.Range("F14").Value = 534.6

# 3) Brand-new rows 15-20. Give the date cells the same "yyyy-mm-dd" number
#    format the existing date column uses so they share that style, before
#    the format code itself is changed to include a time component below.
$ws.Range("A15").Value = 45818
$ws.Range("A15").NumberFormat = "yyyy-mm-dd"
$ws.Range("B15").Value = "SØSTRENE GRENE BOGSTADVEIEN"
$ws.Range("F15").Value = 189

$ws.Range("A16").Value = 45816
$ws.Range("A16").NumberFormat = "yyyy-mm-dd"
$ws.Range("B16").Value = "TILBORDS ALNA"
$ws.Range("F16").Value = 678

$ws.Range("A17").Value = 45814
$ws.Range("A17").NumberFormat = "yyyy-mm-dd"
$ws.Range("B17").Value = "KOMPLETT.NO"
$ws.Range("F17").Value = 2899

$ws.Range("A18").Value = 45813
$ws.Range("A18").NumberFormat = "yyyy-mm-dd"
$ws.Range("B18").Value = "Innbetaling"
$ws.Range("E18").Value = 15000

$ws.Range("A19").Value = 45811
$ws.Range("A19").NumberFormat = "yyyy-mm-dd"
$ws.Range("B19").Value = "COOP OBS STORO"
$ws.Range("F19").Value = 1234.5

$ws.Range("A20").Value = 45809
$ws.Range("A20").NumberFormat = "yyyy-mm-dd"
$ws.Range("B20").Value = "XXL SPORT LAMBERTSETER"
$ws.Range("F20").Value = 1234

# 4) The date column's number format now includes a time component
#    (yyyy-mm-dd -> yyyy-mm-dd h:mm:ss) across every date cell, old and new.
$ws.Range("A2:A20").NumberFormat = "yyyy-mm-dd h:mm:ss"
